$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = -0.02344977690475521
$ws.Cells.Item(2, 4).Value = -0.02401287342951972
$ws.Cells.Item(2, 5).Value = -0.02460367844594259
$ws.Cells.Item(2, 6).Value = -0.02522428873500645
$ws.Cells.Item(2, 7).Value = -0.02587701811145051
$ws.Cells.Item(2, 8).Value = -0.04607133026821583
$ws.Cells.Item(2, 9).Value = -0.04829641013008851
$ws.Cells.Item(2, 10).Value = -0.05074732368792471
$ws.Cells.Item(2, 11).Value = -0.05346029034659393
$ws.Cells.Item(2, 12).Value = -0.05647971215721025
$ws.Cells.Item(2, 13).Value = -0.03398883324524034
$ws.Cells.Item(2, 14).Value = -0.03518472085516693
$ws.Cells.Item(2, 15).Value = -0.03646783132036749
$ws.Cells.Item(2, 16).Value = -0.03784806829059038
$ws.Cells.Item(2, 17).Value = -0.03933689373085549
$ws.Cells.Item(2, 18).Value = -0.0402299362887117
$ws.Cells.Item(2, 19).Value = -0.04191622328075996
$ws.Cells.Item(2, 20).Value = -0.04375006058895344
$ws.Cells.Item(2, 21).Value = -0.04575170024679837
$ws.Cells.Item(2, 22).Value = -0.04794527824532811
$ws.Cells.Item(2, 23).Value = -0.05194558063116002
$ws.Cells.Item(2, 24).Value = -0.05479177098899262
$ws.Cells.Item(2, 25).Value = -0.05796793691303625
$ws.Cells.Item(2, 26).Value = -0.06153499353629193
$ws.Cells.Item(2, 27).Value = -0.06556983277209898
$ws.Cells.Item(2, 28).Value = -0.07141879157394511
$ws.Cells.Item(2, 29).Value = -0.07691173472592662
$ws.Cells.Item(2, 30).Value = -0.08332002216829268
$ws.Cells.Item(2, 31).Value = -0.09089324975262997
$ws.Cells.Item(2, 32).Value = -0.09998083253468056
$ws.Cells.Item(3, 3).Value = -0.02061681392379752
$ws.Cells.Item(3, 4).Value = -0.02105081465243104
$ws.Cells.Item(3, 5).Value = -0.02150348043341723
$ws.Cells.Item(3, 6).Value = -0.02197604181866878
$ws.Cells.Item(3, 7).Value = -0.02246983996131749
$ws.Cells.Item(3, 8).Value = -0.02356563498393421
$ws.Cells.Item(3, 9).Value = -0.02413437690053695
$ws.Cells.Item(3, 10).Value = -0.02473125021443375
$ws.Cells.Item(3, 11).Value = -0.0253583950268798
$ws.Cells.Item(3, 12).Value = -0.02601817416524011
$ws.Cells.Item(3, 13).Value = -0.002871210767738367
$ws.Cells.Item(3, 14).Value = -0.002879478357002462
$ws.Cells.Item(3, 15).Value = -0.002887793696450693
$ws.Cells.Item(3, 16).Value = -0.002896157200959742
$ws.Cells.Item(3, 17).Value = -0.002904569290225315
$ws.Cells.Item(3, 18).Value = -0.0004932023199460452
$ws.Cells.Item(3, 19).Value = -0.000493445688504383
$ws.Cells.Item(3, 20).Value = -0.0004936892973594364
$ws.Cells.Item(3, 21).Value = -0.0004939331468678456
$ws.Cells.Item(3, 22).Value = -0.000494177237385818
$ws.Cells.Item(3, 23).Value = 0.0000922018022508354
$ws.Cells.Item(3, 24).Value = 0.00009219330186224824
$ws.Cells.Item(3, 25).Value = 0.00009218480304068471
$ws.Cells.Item(3, 26).Value = 0.00009217630578628046
$ws.Cells.Item(3, 27).Value = 0.00009216781009803321
$ws.Cells.Item(3, 28).Value = 0.0003759356684309105
$ws.Cells.Item(3, 29).Value = 0.0003757943939142417
$ws.Cells.Item(3, 30).Value = 0.0003756532255380287
$ws.Cells.Item(3, 31).Value = 0.0003755121631827004
$ws.Cells.Item(3, 32).Value = 0.0003753712067288651
$ws.Cells.Item(4, 3).Value = -0.02061681070201453
$ws.Cells.Item(4, 4).Value = -0.02105081129357806
$ws.Cells.Item(4, 5).Value = -0.02150347692855677
$ws.Cells.Item(4, 6).Value = -0.02197603815807006
$ws.Cells.Item(4, 7).Value = -0.02246983613436444
$ws.Cells.Item(4, 8).Value = -0.0235658541973185
$ws.Cells.Item(4, 9).Value = -0.02413460682281439
$ws.Cells.Item(4, 10).Value = -0.02473149164992727
$ws.Cells.Item(4, 11).Value = -0.02535864886252405
$ws.Cells.Item(4, 12).Value = -0.02601844138146701
$ws.Cells.Item(4, 13).Value = -0.04584840409568143
$ws.Cells.Item(4, 14).Value = -0.04805148814138677
$ws.Cells.Item(4, 15).Value = -0.05047698225565728
$ws.Cells.Item(4, 16).Value = -0.05316035663418557
$ws.Cells.Item(4, 17).Value = -0.05614504737593344
$ws.Cells.Item(4, 18).Value = -0.05735934181845381
$ws.Cells.Item(4, 19).Value = -0.06084963694342027
$ws.Cells.Item(4, 20).Value = -0.06479222000764345
$ws.Cells.Item(4, 21).Value = -0.06928109602357364
$ws.Cells.Item(4, 22).Value = -0.07443826028199856
$ws.Cells.Item(4, 23).Value = -0.09274653274091653
$ws.Cells.Item(4, 24).Value = -0.1022278074297299
$ws.Cells.Item(4, 25).Value = -0.1138683156771182
$ws.Cells.Item(4, 26).Value = -0.1285004449018524
$ws.Cells.Item(4, 27).Value = -0.1474475163528693
$ws.Cells.Item(4, 28).Value = -0.1748298019548035
$ws.Cells.Item(4, 29).Value = -0.2118712022913212
$ws.Cells.Item(4, 30).Value = -0.2688281444699049
$ws.Cells.Item(4, 31).Value = -0.3676675222612422
$ws.Cells.Item(4, 32).Value = -0.5814465256885645
$ws.Cells.Item(5, 3).Value = -0.02061559788470009
$ws.Cells.Item(5, 4).Value = -0.02104954687880876
$ws.Cells.Item(5, 5).Value = -0.02150215755219915
$ws.Cells.Item(5, 6).Value = -0.02197466015705228
$ws.Cells.Item(5, 7).Value = -0.02246839551271851
$ws.Cells.Item(5, 8).Value = -0.02356202124482831
$ws.Cells.Item(5, 9).Value = -0.02413058663988751
$ws.Cells.Item(5, 10).Value = -0.0247272701752187
$ws.Cells.Item(5, 11).Value = -0.02535421059057097
$ws.Cells.Item(5, 12).Value = -0.02601376917242309
$ws.Cells.Item(5, 13).Value = -0.03365416220622304
$ws.Cells.Item(5, 14).Value = -0.0348262090961737
$ws.Cells.Item(5, 15).Value = -0.0360828375411654
$ws.Cells.Item(5, 16).Value = -0.03743354610381923
$ws.Cells.Item(5, 17).Value = -0.0388893109169757
$ws.Cells.Item(5, 18).Value = -0.04060736109749505
$ws.Cells.Item(5, 19).Value = -0.04232611284567272
$ws.Cells.Item(5, 20).Value = -0.04419679121819036
$ws.Cells.Item(5, 21).Value = -0.04624047169136423
$ws.Cells.Item(5, 22).Value = -0.04848231689319593
$ws.Cells.Item(5, 23).Value = -0.05835984832431593
$ws.Cells.Item(5, 24).Value = -0.06197680528008735
$ws.Cells.Item(5, 25).Value = -0.06607171936573844
$ws.Cells.Item(5, 26).Value = -0.07074603129146807
$ws.Cells.Item(5, 27).Value = -0.07613207333382718
$ws.Cells.Item(5, 28).Value = -0.08617942711961853
$ws.Cells.Item(5, 29).Value = -0.09430672680959597
$ws.Cells.Item(5, 30).Value = -0.1041265620505165
$ws.Cells.Item(5, 31).Value = -0.1162290984861054
$ws.Cells.Item(5, 32).Value = -0.1315149642141484
$ws.Cells.Item(6, 3).Value = -0.02616037594671301
$ws.Cells.Item(6, 4).Value = -0.02686312540644958
$ws.Cells.Item(6, 5).Value = -0.02760467320454724
$ws.Cells.Item(6, 6).Value = -0.02838832360036028
$ws.Cells.Item(6, 7).Value = -0.0292177670255617
$ws.Cells.Item(6, 8).Value = -0.05011354135718785
$ws.Cells.Item(6, 9).Value = -0.05275740158333174
$ws.Cells.Item(6, 10).Value = -0.05569576544753878
$ws.Cells.Item(6, 11).Value = -0.0589807430800466
$ws.Cells.Item(6, 12).Value = -0.06267750914375136
$ws.Cells.Item(6, 13).Value = -0.04717113088490699
$ws.Cells.Item(6, 14).Value = -0.04950640394503952
$ws.Cells.Item(6, 15).Value = -0.0520849421295595
$ws.Cells.Item(6, 16).Value = -0.05494684539200375
$ws.Cells.Item(6, 17).Value = -0.05814153957805218
$ws.Cells.Item(6, 18).Value = -0.06205188467311738
$ws.Cells.Item(6, 19).Value = -0.06615705459517
$ws.Cells.Item(6, 20).Value = -0.07084387682180356
$ws.Cells.Item(6, 21).Value = -0.07624539628440559
$ws.Cells.Item(6, 22).Value = -0.08253858327495819
$ws.Cells.Item(6, 23).Value = -0.08042473124438186
$ws.Cells.Item(6, 24).Value = -0.08745856263969977
$ws.Cells.Item(6, 25).Value = -0.09584064795204279
$ws.Cells.Item(6, 26).Value = -0.1059997308383307
$ws.Cells.Item(6, 27).Value = -0.1185678958885881
$ws.Cells.Item(6, 28).Value = -0.1718059985881548
$ws.Cells.Item(6, 29).Value = -0.2074465623939228
$ws.Cells.Item(6, 30).Value = -0.2617445746251747
$ws.Cells.Item(6, 31).Value = -0.3545447356411671
$ws.Cells.Item(6, 32).Value = -0.5492940490513409
$ws.Cells.Item(7, 3).Value = -0.02220016115252152
$ws.Cells.Item(7, 4).Value = -0.02270419800711815
$ws.Cells.Item(7, 5).Value = -0.02323165408141564
$ws.Cells.Item(7, 6).Value = -0.02378420039765692
$ws.Cells.Item(7, 7).Value = -0.0243636708270297
$ws.Cells.Item(7, 8).Value = -0.04646504480914875
$ws.Cells.Item(7, 9).Value = -0.04872925167158525
$ws.Cells.Item(7, 10).Value = -0.05122542846735581
$ws.Cells.Item(7, 11).Value = -0.05399114816558226
$ws.Cells.Item(7, 12).Value = -0.05707256127771673
$ws.Cells.Item(7, 13).Value = -0.01096836947073912
$ws.Cells.Item(7, 14).Value = -0.01109000878452149
$ws.Cells.Item(7, 15).Value = -0.01121437631638311
$ws.Cells.Item(7, 16).Value = -0.01134156489311114
$ws.Cells.Item(7, 17).Value = -0.01147167160100641
$ws.Cells.Item(7, 18).Value = -0.007514405846642954
$ws.Cells.Item(7, 19).Value = -0.007571299665113281
$ws.Cells.Item(7, 20).Value = -0.007629061576472239
$ws.Cells.Item(7, 21).Value = -0.007687711601663607
$ws.Cells.Item(7, 22).Value = -0.007747270382060884
$ws.Cells.Item(7, 23).Value = -0.006996024207411306
$ws.Cells.Item(7, 24).Value = -0.007045313390439621
$ws.Cells.Item(7, 25).Value = -0.007095302016747304
$ws.Cells.Item(7, 26).Value = -0.007146005081010405
$ws.Cells.Item(7, 27).Value = -0.007197438009597244
$ws.Cells.Item(7, 28).Value = -0.006613902959854232
$ws.Cells.Item(7, 29).Value = -0.006657937915137687
$ws.Cells.Item(7, 30).Value = -0.006702563164559738
$ws.Cells.Item(7, 31).Value = -0.006747790657663831
$ws.Cells.Item(7, 32).Value = -0.006793632668717403
$ws.Cells.Item(8, 3).Value = -0.02602168995290594
$ws.Cells.Item(8, 4).Value = -0.0267169090774185
$ws.Cells.Item(8, 5).Value = -0.02745029614363624
$ws.Cells.Item(8, 6).Value = -0.02822508303153047
$ws.Cells.Item(8, 7).Value = -0.02904487709929878
$ws.Cells.Item(8, 8).Value = -0.04999675967651204
$ws.Cells.Item(8, 9).Value = -0.05262798857347846
$ws.Cells.Item(8, 10).Value = -0.05555155518499324
$ws.Cells.Item(8, 11).Value = -0.05881904458625517
$ws.Cells.Item(8, 12).Value = -0.06249493707657749
$ws.Cells.Item(8, 13).Value = -0.04707652434093104
$ws.Cells.Item(8, 14).Value = -0.04940220861740389
$ws.Cells.Item(8, 15).Value = -0.05196962276290482
$ws.Cells.Item(8, 16).Value = -0.05481852059884745
$ws.Cells.Item(8, 17).Value = -0.05799787849586219
$ws.Cells.Item(8, 18).Value = -0.06137467986040702
$ws.Cells.Item(8, 19).Value = -0.06538783745070859
$ws.Cells.Item(8, 20).Value = -0.06996253640906373
$ws.Cells.Item(8, 21).Value = -0.07522550343180127
$ws.Cells.Item(8, 22).Value = -0.08134469939532299
$ws.Cells.Item(8, 23).Value = -0.08860758695431899
$ws.Cells.Item(8, 24).Value = -0.09722221261225024
$ws.Cells.Item(8, 25).Value = -0.1076922959010428
$ws.Cells.Item(8, 26).Value = -0.1206896403632302
$ws.Cells.Item(8, 27).Value = -0.1372548828073461
$ws.Cells.Item(8, 28).Value = -0.1697085398284158
$ws.Cells.Item(8, 29).Value = -0.2043963451019292
$ws.Cells.Item(8, 30).Value = -0.2569072475265534
$ws.Cells.Item(8, 31).Value = -0.3457270262311347
$ws.Cells.Item(8, 32).Value = -0.5284140413742193
